$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "PDL Date" values (01012025) with 99999999 for all data rows (B2:B267)
$ws.Range("B2:B267").Value = "99999999"
